# Refresh the Spearman-correlation table (keywords + correlation values) in Sheet1.
# The underlying dataset changed: keyword list (col B) was replaced/reordered, the
# correlation values (cols C:F) were recomputed, and six new keyword rows were appended
# (rows 54-59), extending the used range from A1:F53 to A1:F59.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2=0)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "reabertura comercio"
$ws.Range("C2").Value = 0.5410279032331542
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
# Row 3 (A3=1)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "reabertura shopping"
$ws.Range("C3").Value = 0.6621219850108133
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
# Row 4 (A4=2)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "reabertura escolas"
$ws.Range("C4").Value = 0.6997480710751741
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
# Row 5 (A5=3)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "teste rapido de coronavírus"
$ws.Range("C5").Value = 0.585868109602798
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
# Row 6 (A6=4)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "oximetro"
$ws.Range("C6").Value = 0.5881031926710265
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
# Row 7 (A7=5)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "covid"
$ws.Range("C7").Value = 0.5006347415728153
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
# Row 8 (A8=6)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "covid pcr"
$ws.Range("C8").Value = 0.5844180058725504
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").Value = 0.5844275039268546
# Row 9 (A9=7)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "pcr exame covid"
$ws.Range("C9").Value = 0.6019890928692534
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").Value = 0.6567849052317277
# Row 10 (A10=8)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "pcr"
$ws.Range("C10").Value = 0.6052549139284656
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").Value = 0.6670784283315848
# Row 11 (A11=9)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "decreto lockdown"
$ws.Range("C11").Value = 0.6360316878370555
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
# Row 12 (A12=10)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "estou com covid"
$ws.Range("C12").Value = 0.6579449281943737
$ws.Range("D12").ClearContents()
$ws.Range("E12").Value = 0.6365432545821372
$ws.Range("F12").ClearContents()
# Row 13 (A13=11)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "peguei covid"
$ws.Range("C13").Value = 0.6816659979672126
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").Value = 0.6413749437967033
# Row 14 (A14=12)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "febre"
$ws.Range("C14").ClearContents()
$ws.Range("D14").Value = -0.5204727820886367
$ws.Range("E14").Value = -0.6641040680504862
$ws.Range("F14").ClearContents()
# Row 15 (A15=13)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "coronavírus no brasil"
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = -0.5481748671481036
$ws.Range("E15").Value = -0.6924857097767528
$ws.Range("F15").Value = -0.6470189149802471
# Row 16 (A16=14)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "ministerio da saude"
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = -0.5597582214328947
$ws.Range("E16").Value = -0.6957502478314564
$ws.Range("F16").Value = -0.6452403588393268
# Row 17 (A17=15)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "taxa de ocupação de leitos"
$ws.Range("C17").ClearContents()
$ws.Range("D17").Value = 0.5964688143240882
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()
# Row 18 (A18=16)
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "teste rápido covid"
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = 0.6201524044070241
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()
# Row 19 (A19=17)
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "teste igg"
$ws.Range("C19").ClearContents()
$ws.Range("D19").Value = 0.6340277362477
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
# Row 20 (A20=18)
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "mortes corona"
$ws.Range("C20").ClearContents()
$ws.Range("D20").Value = -0.5283548259416021
$ws.Range("E20").Value = -0.6635859979449293
$ws.Range("F20").Value = -0.6247949728432902
# Row 21 (A21=19)
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "oms corona"
$ws.Range("C21").ClearContents()
$ws.Range("D21").Value = -0.6065706654261432
$ws.Range("E21").ClearContents()
$ws.Range("F21").ClearContents()
# Row 22 (A22=20)
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "álcool gel"
$ws.Range("C22").ClearContents()
$ws.Range("D22").Value = -0.5818496783484663
$ws.Range("E22").Value = -0.6983307520119921
$ws.Range("F22").Value = -0.6938699973773165
# Row 23 (A23=21)
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "álcool 70"
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = -0.6300459277469465
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()
# Row 24 (A24=22)
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "medicamento coronavírus"
$ws.Range("C24").ClearContents()
$ws.Range("D24").Value = -0.5322463163450158
$ws.Range("E24").Value = -0.6505220596301851
$ws.Range("F24").Value = -0.6093253982597782
# Row 25 (A25=23)
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "medicamento corona"
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = -0.5783872290466793
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()
# Row 26 (A26=24)
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "teste covid"
$ws.Range("C26").ClearContents()
$ws.Range("D26").ClearContents()
$ws.Range("E26").Value = 0.6186716308546935
$ws.Range("F26").ClearContents()
# Row 27 (A27=25)
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "teste coronavírus"
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E27").Value = -0.5115036355418436
$ws.Range("F27").Value = -0.5437209309855231
# Row 28 (A28=26)
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "tosse"
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("E28").Value = -0.5797123425191882
$ws.Range("F28").ClearContents()
# Row 29 (A29=27)
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "sintomas coronavirus"
$ws.Range("C29").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("E29").Value = -0.6575612014236544
$ws.Range("F29").ClearContents()
# Row 30 (A30=28)
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "ivermectina covid como tomar"
$ws.Range("C30").ClearContents()
$ws.Range("D30").ClearContents()
$ws.Range("E30").Value = 0.651538034628236
$ws.Range("F30").Value = 0.5575742512017
# Row 31 (A31=29)
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "pico brasil"
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").Value = -0.5206594072721696
$ws.Range("F31").ClearContents()
# Row 32 (A32=30)
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "pico coronavírus"
$ws.Range("C32").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("E32").Value = -0.5061291431008383
$ws.Range("F32").ClearContents()
# Row 33 (A33=31)
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "brasil coronavírus"
$ws.Range("C33").ClearContents()
$ws.Range("D33").ClearContents()
$ws.Range("E33").Value = -0.6500685980652502
$ws.Range("F33").Value = -0.5995637481001683
# Row 34 (A34=32)
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "coronavírus brasil pico"
$ws.Range("C34").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("E34").Value = -0.5515835213760795
$ws.Range("F34").Value = -0.5296711575521518
# Row 35 (A35=33)
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "coronavírus pico brasil"
$ws.Range("C35").ClearContents()
$ws.Range("D35").ClearContents()
$ws.Range("E35").Value = -0.5515835213760795
$ws.Range("F35").Value = -0.5296711575521518
# Row 36 (A36=34)
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "pico corona"
$ws.Range("C36").ClearContents()
$ws.Range("D36").ClearContents()
$ws.Range("E36").Value = -0.5137387578764975
$ws.Range("F36").Value = -0.5618329471787832
# Row 37 (A37=35)
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "corona grupo de risco"
$ws.Range("C37").ClearContents()
$ws.Range("D37").ClearContents()
$ws.Range("E37").Value = -0.5810904045939116
$ws.Range("F37").Value = -0.5889948687900587
# Row 38 (A38=36)
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "leitos uti brasil"
$ws.Range("C38").ClearContents()
$ws.Range("D38").ClearContents()
$ws.Range("E38").Value = -0.5196117308810991
$ws.Range("F38").Value = -0.573959540188101
# Row 39 (A39=37)
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "máscara n95"
$ws.Range("C39").ClearContents()
$ws.Range("D39").ClearContents()
$ws.Range("E39").Value = -0.5826930393282175
$ws.Range("F39").Value = -0.615803630710012
# Row 40 (A40=38)
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "máscara descartável"
$ws.Range("C40").ClearContents()
$ws.Range("D40").ClearContents()
$ws.Range("E40").Value = -0.5759033019702267
$ws.Range("F40").Value = -0.6438451327522168
# Row 41 (A41=39)
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "igg"
$ws.Range("C41").ClearContents()
$ws.Range("D41").ClearContents()
$ws.Range("E41").Value = 0.5809379031725802
$ws.Range("F41").ClearContents()
# Row 42 (A42=40)
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "igm"
$ws.Range("C42").ClearContents()
$ws.Range("D42").ClearContents()
$ws.Range("E42").Value = 0.6073531272539124
$ws.Range("F42").ClearContents()
# Row 43 (A43=41)
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "igg igm"
$ws.Range("C43").ClearContents()
$ws.Range("D43").ClearContents()
$ws.Range("E43").Value = 0.5903755247221089
$ws.Range("F43").Value = 0.500688086829228
# Row 44 (A44=42)
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "covid igg"
$ws.Range("C44").ClearContents()
$ws.Range("D44").ClearContents()
$ws.Range("E44").Value = 0.571398530421381
$ws.Range("F44").ClearContents()
# Row 45 (A45=43)
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "covid igm igg"
$ws.Range("C45").ClearContents()
$ws.Range("D45").ClearContents()
$ws.Range("E45").Value = 0.5734232018632055
$ws.Range("F45").ClearContents()
# Row 46 (A46=44)
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "exame igg"
$ws.Range("C46").ClearContents()
$ws.Range("D46").ClearContents()
$ws.Range("E46").Value = 0.608752531678521
$ws.Range("F46").ClearContents()
# Row 47 (A47=45)
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "coronavírus quarentena"
$ws.Range("C47").ClearContents()
$ws.Range("D47").ClearContents()
$ws.Range("E47").Value = -0.654865208061991
$ws.Range("F47").Value = -0.6312004276824037
# Row 48 (A48=46)
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "coronavírus idosos"
$ws.Range("C48").ClearContents()
$ws.Range("D48").ClearContents()
$ws.Range("E48").Value = -0.5910193844118176
$ws.Range("F48").Value = -0.613621007743175
# Row 49 (A49=47)
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = "coronavírus oms"
$ws.Range("C49").ClearContents()
$ws.Range("D49").ClearContents()
$ws.Range("E49").Value = -0.5648388009144159
$ws.Range("F49").Value = -0.5130695531360026
# Row 50 (A50=48)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "oms coronavírus"
$ws.Range("C50").ClearContents()
$ws.Range("D50").ClearContents()
$ws.Range("E50").Value = -0.5648388009144159
$ws.Range("F50").Value = -0.5130695531360026
# Row 51 (A51=49)
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "oms brasil coronavírus"
$ws.Range("C51").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("E51").Value = -0.5003870462679144
$ws.Range("F51").ClearContents()
# Row 52 (A52=50)
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "coronavírus imunidade"
$ws.Range("C52").ClearContents()
$ws.Range("D52").ClearContents()
$ws.Range("E52").Value = -0.5107978788174282
$ws.Range("F52").Value = -0.5067484612367316
# Row 53 (A53=51)
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = "respirador mecanico"
$ws.Range("C53").ClearContents()
$ws.Range("D53").ClearContents()
$ws.Range("E53").Value = -0.6265050702644759
$ws.Range("F53").Value = -0.660703172827806
# Row 54 (A54=52)
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "respirador"
$ws.Range("C54").ClearContents()
$ws.Range("D54").ClearContents()
$ws.Range("E54").Value = -0.5156948738712612
$ws.Range("F54").Value = -0.5152082519872652
# Row 55 (A55=53)
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "dor de garganta"
$ws.Range("C55").ClearContents()
$ws.Range("D55").ClearContents()
$ws.Range("E55").ClearContents()
$ws.Range("F55").Value = -0.5971462949045846
# Row 56 (A56=54)
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = "exame cotonete"
$ws.Range("C56").ClearContents()
$ws.Range("D56").ClearContents()
$ws.Range("E56").ClearContents()
$ws.Range("F56").Value = 0.6686192510512416
# Row 57 (A57=55)
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = "tomar ivermectina"
$ws.Range("C57").ClearContents()
$ws.Range("D57").ClearContents()
$ws.Range("E57").ClearContents()
$ws.Range("F57").Value = 0.6318613438932613
# Row 58 (A58=56)
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = "exame covid"
$ws.Range("C58").ClearContents()
$ws.Range("D58").ClearContents()
$ws.Range("E58").ClearContents()
$ws.Range("F58").Value = 0.5562220533565535
# Row 59 (A59=57)
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = "reagente igg"
$ws.Range("C59").ClearContents()
$ws.Range("D59").ClearContents()
$ws.Range("E59").ClearContents()
$ws.Range("F59").Value = 0.6963656043963102
